$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C19:F19").Value = 5

# Row 19's C:F cells were style "s=5" (green highlighted, thick border).
# The target state uses style "s=2" (same thick border, no fill) - the
# same style already used by e.g. C21:F21. Copy formats from there so we
# reuse the existing style index instead of minting a new one.
$ws.Range("C21:F21").Copy()
$ws.Range("C19:F19").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C19:F19").Select()
